$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 3-8 (columns D..T) get cyclically reshuffled.
# Mapping: new row N gets the old D..T values that used to live in row M:
#   3 <- 5, 4 <- 3, 5 <- 7, 6 <- 8, 7 <- 4, 8 <- 6

$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Capture the original values for rows 3..8 before overwriting anything.
$orig = @{}
foreach ($r in 3..8) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $orig[$r] = $rowVals
}

$mapping = @{ 3 = 5; 4 = 3; 5 = 7; 6 = 8; 7 = 4; 8 = 6 }

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
